$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 22631.215
$ws.Range("I62").Value = 35263.43
$ws.Range("K62").Value = 35263.43
$ws.Range("M62").Value = -34639.43

$ws.Range("H65").Value = 22631.215
$ws.Range("I65").Value = 35263.43
$ws.Range("K65").Value = 176317.15
$ws.Range("M65").Value = -173197.15

$ws.Range("H106").Value = 3626.25
$ws.Range("I106").Value = 3626.25
$ws.Range("K106").Value = 3626.25
$ws.Range("M106").Value = -2995.25

$ws.Range("H107").Value = 1567.7333
$ws.Range("I107").Value = 899.2381
$ws.Range("J107").Value = 3127.5557
$ws.Range("K107").Value = 899.2381
$ws.Range("L107").Value = 3127.5557
$ws.Range("M107").Value = 1020.7619
$ws.Range("N107").Value = -6967.5557

$ws.Range("H116").Value = 6949.067
$ws.Range("I116").Value = 6841.1816
$ws.Range("J116").Value = 7245.75
$ws.Range("K116").Value = 6841.1816
$ws.Range("L116").Value = 7245.75
$ws.Range("M116").Value = -3399.1816
$ws.Range("N116").Value = -14129.75

$ws.Range("H132").Value = 3550
$ws.Range("I132").Value = 2574.8372
$ws.Range("J132").Value = 11936.4
$ws.Range("K132").Value = 7724.5116
$ws.Range("L132").Value = 35809.2
$ws.Range("M132").Value = -5194.5116
$ws.Range("N132").Value = -40869.2

$ws.Range("H138").Value = 3513.4746
$ws.Range("I138").Value = 3124
$ws.Range("J138").Value = 3602.7292
$ws.Range("K138").Value = 9372
$ws.Range("L138").Value = 10808.1876
$ws.Range("M138").Value = -4232
$ws.Range("N138").Value = -21088.1876

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1926.6428
$ws.Range("J94").Value = 4033
$ws.Range("L94").Value = 4033
$ws.Range("N94").Value = -4935

$ws.Range("H99").Value = 2076.9
$ws.Range("I99").Value = 1619.9412
$ws.Range("K99").Value = 1619.9412
$ws.Range("M99").Value = -121.9412

$ws.Range("H105").Value = 2854.2307
$ws.Range("I105").Value = 2345.3333
$ws.Range("K105").Value = 2345.3333
$ws.Range("M105").Value = -598.3332999999998

$ws.Range("H107").Value = 1063.0588
$ws.Range("I107").Value = 911.4666999999999
$ws.Range("K107").Value = 911.4666999999999
$ws.Range("M107").Value = 1008.5333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1191.6666
$ws.Range("I16").Value = 716
$ws.Range("K16").Value = 716
$ws.Range("M16").Value = -429

$ws.Range("H31").Value = 2660.3809
$ws.Range("I31").Value = 1304.6129
$ws.Range("J31").Value = 6481.1816
$ws.Range("K31").Value = 1304.6129
$ws.Range("L31").Value = 6481.1816
$ws.Range("M31").Value = -1009.6129
$ws.Range("N31").Value = -7071.1816

$ws.Range("H34").Value = 2660.3809
$ws.Range("I34").Value = 1304.6129
$ws.Range("J34").Value = 6481.1816
$ws.Range("K34").Value = 1304.6129
$ws.Range("L34").Value = 6481.1816
$ws.Range("M34").Value = -1102.6129
$ws.Range("N34").Value = -6885.1816

$ws.Range("H86").Value = 37395.418
$ws.Range("I86").Value = 44075.6
$ws.Range("K86").Value = 44075.6
$ws.Range("M86").Value = -42952.6

$ws.Range("H89").Value = 37395.418
$ws.Range("I89").Value = 44075.6
$ws.Range("K89").Value = 220378
$ws.Range("M89").Value = -214762

$ws.Range("H105").Value = 5964.2
$ws.Range("I105").Value = 3940.6667
$ws.Range("K105").Value = 3940.6667
$ws.Range("M105").Value = -2193.6667

$ws.Range("H107").Value = 397.4
$ws.Range("I107").Value = 397.4
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 397.4
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1522.6
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 1191.6666
$ws.Range("I113").Value = 716
$ws.Range("K113").Value = 716
$ws.Range("M113").Value = 1454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2555008.5
$ws.Range("I4").Value = 2122566
$ws.Range("K4").Value = 6367698
$ws.Range("M4").Value = -6367586

$ws.Range("H5").Value = 1150.9333
$ws.Range("I5").Value = 1559.3334
$ws.Range("K5").Value = 4678.0002
$ws.Range("M5").Value = -4566.0002

$ws.Range("H68").Value = 763.8
$ws.Range("I68").Value = 766.44446
$ws.Range("J68").Value = 740
$ws.Range("K68").Value = 2299.33338
$ws.Range("L68").Value = 2220
$ws.Range("M68").Value = -1488.33338
$ws.Range("N68").Value = -3842

$ws.Range("H71").Value = 763.8
$ws.Range("I71").Value = 766.44446
$ws.Range("J71").Value = 740
$ws.Range("K71").Value = 6898.00014
$ws.Range("L71").Value = 6660
$ws.Range("M71").Value = -2842.00014
$ws.Range("N71").Value = -14772

$ws.Range("H122").Value = 486.45834
$ws.Range("J122").Value = 588.36365
$ws.Range("L122").Value = 5295.27285
$ws.Range("N122").Value = -10195.27285

$ws.Range("H135").Value = 1150.9333
$ws.Range("I135").Value = 1559.3334
$ws.Range("K135").Value = 14034.0006
$ws.Range("M135").Value = -11499.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6554.6665
$ws.Range("I70").Value = 6000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5730

$ws.Range("H73").Value = 6554.6665
$ws.Range("I73").Value = 6000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -5064

$ws.Range("H97").Value = 724.1053000000001
$ws.Range("I97").Value = 724.1053000000001
$ws.Range("K97").Value = 724.1053000000001
$ws.Range("M97").Value = -228.1053000000001

$ws.Range("H107").Value = 78025.69500000001
$ws.Range("I107").Value = 143120
$ws.Range("J107").Value = 2082.3333
$ws.Range("K107").Value = 143120
$ws.Range("L107").Value = 2082.3333
$ws.Range("M107").Value = -141200
$ws.Range("N107").Value = -5922.3333

$ws.Range("H113").Value = 251722.12
$ws.Range("I113").Value = 334602.66
$ws.Range("J113").Value = 201993.8
$ws.Range("K113").Value = 334602.66
$ws.Range("L113").Value = 201993.8
$ws.Range("M113").Value = -332432.66
$ws.Range("N113").Value = -206333.8

$ws.Range("H132").Value = 33175.637
$ws.Range("I132").Value = 47599.727
$ws.Range("K132").Value = 142799.181
$ws.Range("M132").Value = -140269.181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8928.607
$ws.Range("I7").Value = 10565.55
$ws.Range("K7").Value = 10565.55
$ws.Range("M7").Value = -10453.55

$ws.Range("H40").Value = 2631.889
$ws.Range("I40").Value = 1549.2142
$ws.Range("J40").Value = 6421.25
$ws.Range("K40").Value = 1549.2142
$ws.Range("L40").Value = 6421.25
$ws.Range("M40").Value = -1413.2142
$ws.Range("N40").Value = -6693.25

$ws.Range("H55").Value = 2416.25
$ws.Range("I55").Value = 2429.1
$ws.Range("K55").Value = 2429.1
$ws.Range("M55").Value = -2256.1

$ws.Range("H122").Value = 4085.3489
$ws.Range("I122").Value = 3440
$ws.Range("J122").Value = 4900.5264
$ws.Range("K122").Value = 10320
$ws.Range("L122").Value = 14701.5792
$ws.Range("M122").Value = -7870
$ws.Range("N122").Value = -19601.5792

$ws.Range("H126").Value = 8928.607
$ws.Range("I126").Value = 10565.55
$ws.Range("K126").Value = 31696.65
$ws.Range("M126").Value = -29226.65

$ws.Range("H132").Value = 35462.082
$ws.Range("I132").Value = 41261.195
$ws.Range("K132").Value = 123783.585
$ws.Range("M132").Value = -121253.585

$ws.Range("H133").Value = 72000
$ws.Range("J133").Value = 72000
$ws.Range("L133").Value = 72000
$ws.Range("N133").Value = -77060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 637.9167
$ws.Range("J107").Value = 1528
$ws.Range("L107").Value = 4584
$ws.Range("N107").Value = -8424

$ws.Range("H132").Value = 46979.348
$ws.Range("I132").Value = 53564.75
$ws.Range("J132").Value = 3076.6667
$ws.Range("K132").Value = 160694.25
$ws.Range("L132").Value = 9230.000100000001
$ws.Range("M132").Value = -158164.25
$ws.Range("N132").Value = -14290.0001

$ws.Range("H133").Value = 84000
$ws.Range("J133").Value = 84000
$ws.Range("L133").Value = 84000
$ws.Range("N133").Value = -94120
